$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column F (6) width wider to fit long fastq filenames
$ws.Columns.Item(6).ColumnWidth = 49.86

# Update manualStatus column (I) from raw numeric codes to bracketed text labels
$ws.Range("I5").Value = "[4]"
$ws.Range("I8").Value = "[6]"
$ws.Range("I9").Value = "[4]"
$ws.Range("I10").Value = "[6]"
$ws.Range("I19").Value = "[4]"
$ws.Range("I20").Value = "[4]"
$ws.Range("I21").Value = "[4]"
$ws.Range("I22").Value = "[6]"
$ws.Range("I23").Value = "[4]"
$ws.Range("I24").Value = "[4]"

# Row heights settle slightly shorter once the manual status text is entered
$ws.Rows.Item(9).RowHeight = 13.8
$ws.Rows.Item(10).RowHeight = 13.8
$ws.Rows.Item(19).RowHeight = 13.8
$ws.Rows.Item(20).RowHeight = 13.8
$ws.Rows.Item(21).RowHeight = 13.8
$ws.Rows.Item(22).RowHeight = 13.8
$ws.Rows.Item(23).RowHeight = 13.8
$ws.Rows.Item(24).RowHeight = 13.8

# Reflect last-edited selection
$ws.Range("I22").Select()
